# TFS 26411 - New coaching reason for Claims View (Medicare Only)
# Adds rows to DIM_Coaching_Reason, DIM_Sub_Coaching_Reason, Coaching_Reason_Selection
# and logs the change on Revision_History.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) Revision_History - log the change, un-highlight the previous latest row
# ---------------------------------------------------------------------------
$wsRev = $wb.Worksheets.Item("Revision_History")

# Copy the current (highlighted) formatting of row 117 onto the new row 118
# BEFORE we change row 117's own formatting.
$wsRev.Range("A117:E117").Copy()
$wsRev.Range("A118:E118").PasteSpecial($xlPasteFormats)

# Un-highlight row 117 (it is no longer the most recent entry) by copying the
# formatting from row 116, which already uses the "older entry" style.
$wsRev.Range("A116:E116").Copy()
$wsRev.Range("A117:E117").PasteSpecial($xlPasteFormats)

# New row 118 values
$wsRev.Range("A118").Value2 = 107
$wsRev.Range("B118").Value2 = 45019
$wsRev.Range("C118").Value2 = "Susmitha Palacherla"
$wsRev.Range("D118").Value2 = 26411
$wsRev.Range("E118").Value2 = "New coaching reason for Claims View (Medicare Only).Added row(s) to DIM_Coaching_Reason, DIM_Sub_Coaching_Reason and Coaching_Reason_Selection"

$wsRev.Range("A118:XFD118").Select()

# ---------------------------------------------------------------------------
# 2) DIM_Coaching_Reason - add new coaching reason
# ---------------------------------------------------------------------------
$wsDCR = $wb.Worksheets.Item("DIM_Coaching_Reason")

$wsDCR.Range("A75:B75").Interior.Color = 65535
$wsDCR.Range("A75").Value2 = 73
$wsDCR.Range("B75").Value2 = "Claims View (Medicare Only)"

$wsDCR.Range("A75:B75").Select()

# ---------------------------------------------------------------------------
# 3) DIM_Sub_Coaching_Reason - add new sub coaching reasons
# ---------------------------------------------------------------------------
$wsDSCR = $wb.Worksheets.Item("DIM_Sub_Coaching_Reason")

$subReasons = @(
  "Old View",
  "Claim Header",
  "Claim Summary",
  "Crossover",
  "Patient Responsibility",
  "Preventive Service Pop-Up",
  "Benefit Period Pop-Up",
  "Check Research",
  "Provider Information",
  "Undeliverable Address indicator",
  "Claim Status"
)

$row = 305
$id = 303
foreach ($reason in $subReasons) {
  $wsDSCR.Range("A" + $row + ":B" + $row).Interior.Color = 65535
  $wsDSCR.Range("A" + $row).Value2 = $id
  $wsDSCR.Range("B" + $row).Value2 = $reason
  $row = $row + 1
  $id = $id + 1
}

$wsDSCR.Range("A305:B315").Select()

# ---------------------------------------------------------------------------
# 4) DIM_Site / Survey_Sites - un-highlight rows from the previous change
# ---------------------------------------------------------------------------
$wsSite = $wb.Worksheets.Item("DIM_Site")
$wsSite.Range("A19:E19").Copy()
$wsSite.Range("A20:E20").PasteSpecial($xlPasteFormats)
$wsSite.Range("A20:E20").Select()

$wsSurveySites = $wb.Worksheets.Item("Survey_Sites")
$wsSurveySites.Range("A19:E19").Copy()
$wsSurveySites.Range("A20:E20").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# 5) Coaching_Reason_Selection - wire the new coaching/sub-coaching reasons
#    (kept last so this ends up the active sheet/tab, matching the workbook)
# ---------------------------------------------------------------------------
$wsCRS = $wb.Worksheets.Item("Coaching_Reason_Selection")

$subIds = @(303,304,305,306,307,308,309,310,311,312,313)
$subTexts = $subReasons

$row = 446
for ($i = 0; $i -lt $subIds.Length; $i++) {
  $rng = $wsCRS.Range("A" + $row + ":P" + $row)
  $rng.Interior.Color = 65535
  $wsCRS.Range("A" + $row).Value2 = 73
  $wsCRS.Range("B" + $row).Value2 = "Claims View (Medicare Only)"
  $wsCRS.Range("C" + $row).Value2 = $subIds[$i]
  $wsCRS.Range("D" + $row).Value2 = $subTexts[$i]
  $wsCRS.Range("E" + $row + ":J" + $row).Value2 = 1
  $wsCRS.Range("K" + $row + ":P" + $row).Value2 = 0
  $row = $row + 1
}

# Final row (457) reuses the existing "Other" sub coaching reason (id 42)
$rng = $wsCRS.Range("A457:P457")
$rng.Interior.Color = 65535
$wsCRS.Range("A457").Value2 = 73
$wsCRS.Range("B457").Value2 = "Claims View (Medicare Only)"
$wsCRS.Range("C457").Value2 = 42
$wsCRS.Range("D457").Value2 = "Other: Specify reason under coaching details."
$wsCRS.Range("E457:J457").Value2 = 1
$wsCRS.Range("K457:P457").Value2 = 0

$wsCRS.Activate()
$wsCRS.Range("A446:P457").Select()

Write-Output "edits applied"
